$wb = $excel.ActiveWorkbook
$wsMaterials = $wb.Worksheets.Item("Materials")
$wsRoute = $wb.Worksheets.Item("Route 1")

# --- 1) Update comment text on Route 1 D1 (Volumes) ---
$d1Comment = "Ryan Nelson:`nFor solvents, express the amount as volumes relative to a certain compound. E.g 1 volume = 1 L solvent per 1 kg of SM. (Which is equal to 1 mL per 1 g, etc.) If a `"Volume`" is given, then `"Density`" of the compound must be provided in the Materials tab."
$wsRoute.Range("D1").Comment.Text($d1Comment)

# --- 2) Update comment text on Route 1 E1 (Sol Recyc -> Recycle) ---
$e1Comment = "Ryan Nelson:`n(Optional) The fractional amount of Compound that can be expected to be recycled/recovered. E.g. `"0`" indicates that none of this compound can be recycled; `"1`" indicates that 100% of this compound can be recycled."
$wsRoute.Range("E1").Comment.Text($e1Comment)

# --- 3) Capture the Notes comment (currently on H1) before the columns move ---
$notesComment = $wsRoute.Range("H1").Comment.Text()

# --- 4) Re-point the F1 comment (currently "Cost step") to the Notes text, since
#        after the column deletion below F1 becomes the old "Notes" column. Reusing
#        the existing comment object (instead of delete+AddComment) keeps the
#        original "Ryan Nelson" authorship instead of a new generic author. ---
$wsRoute.Range("F1").Comment.Text($notesComment)

# --- 5) Remove the comments that belonged to the columns being deleted (OPEX / the old Notes spot) ---
$wsRoute.Range("G1").Comment.Delete()
$wsRoute.Range("H1").Comment.Delete()

# --- 6) Delete the "Cost step" (F) and "OPEX" (G) columns entirely (rightmost first) ---
$wsRoute.Columns.Item(7).Delete()
$wsRoute.Columns.Item(6).Delete()

# --- 7) Rename the E1 header label from "Sol Recyc" to "Recycle" ---
$wsRoute.Range("E1").Value = "Recycle"

# --- 8) Tidy up the conditional formatting so it matches the current columns ---
# The "whole row" highlight rule used to cover A2:E1048576; narrow it back down
# to the A column now that each column manages its own pair of rules.
$fcA2 = $wsRoute.Range("A2").FormatConditions.Item(2)
$fcA2.ModifyAppliesToRange($wsRoute.Range("A2:A1048576"))

# Column B now gets the same two-rule pattern already used elsewhere (Materials!B).
$fcB2 = $wsRoute.Range("B2:B1048576").FormatConditions.Add(2, 0, "=NOT(ISBLANK(B2))")
$fcB2.Interior.Color = 13431551

# The old "Sol Recyc"/E-column rule tied to the now-removed layout no longer applies.
$wsRoute.Range("E2").FormatConditions.Item(1).Delete()

# --- 9) Restore the selections/active cells as seen in the saved workbook ---
$wsRoute.Activate()
$wsRoute.Range("D15").Select()
$wsMaterials.Activate()
$wsMaterials.Range("A2").Select()
